$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 182
    3  = 185
    4  = 187
    5  = 189
    6  = 190
    7  = 193
    8  = 194
    9  = 196
    10 = 198
    11 = 200
    12 = 202
    13 = 203
    14 = 205
    15 = 207
    16 = 32
    17 = 126
    18 = 143
    19 = 157
    20 = 210
    21 = 244
    22 = 290
    23 = 381
    24 = 430
    25 = 461
    26 = 492
    27 = 503
}

foreach ($row in $values.Keys) {
    $ws.Range("A$row").Value = $values[$row]
}
